$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report number and week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/9/2024  Through  12/15/2024"

# --- Crime Complaints table updates (rows 15-33) ---
$ws.Range("N15").Value = -33.333333333333
$ws.Range("C16").Value = 1
$ws.Range("D16").NumberFormat = '#,##0'
$ws.Range("D16").Value = 3
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 28.571428571428
$ws.Range("I16").Value = 113
$ws.Range("J16").Value = 136
$ws.Range("K16").Value = -16.911764705882
$ws.Range("L16").Value = -37.56906077348
$ws.Range("M16").Value = -54.618473895582
$ws.Range("N16").Value = -88.574317492416
$ws.Range("C17").NumberFormat = '#,##0'
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -55
$ws.Range("I17").Value = 262
$ws.Range("J17").Value = 238
$ws.Range("K17").Value = 10.084033613445
$ws.Range("L17").Value = 29.064039408867
$ws.Range("M17").Value = 167.34693877551
$ws.Range("N17").Value = 14.410480349345
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 75
$ws.Range("F18").Value = 29
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 38.095238095238
$ws.Range("I18").Value = 221
$ws.Range("J18").Value = 255
$ws.Range("K18").Value = -13.333333333333
$ws.Range("L18").Value = -2.21238938053
$ws.Range("M18").Value = -18.148148148148
$ws.Range("N18").Value = -83.821376281112
$ws.Range("C19").Value = 11
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -17.777777777777
$ws.Range("I19").Value = 520
$ws.Range("J19").Value = 603
$ws.Range("K19").Value = -13.764510779436
$ws.Range("L19").Value = -17.721518987341
$ws.Range("M19").Value = 19.540229885057
$ws.Range("N19").Value = -10.958904109589
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -33.333333333333
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 81.25
$ws.Range("I20").Value = 369
$ws.Range("J20").Value = 329
$ws.Range("K20").Value = 12.158054711246
$ws.Range("L20").Value = 40.304182509505
$ws.Range("M20").Value = 66.216216216216
$ws.Range("N20").Value = -92.213547161848
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -6.666666666666
$ws.Range("F21").Value = 113
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = 3.669724770642
$ws.Range("I21").Value = 1509
$ws.Range("J21").Value = 1576
$ws.Range("K21").Value = -4.251269035532
$ws.Range("L21").Value = -1.114023591087
$ws.Range("M21").Value = 16.88613477924
$ws.Range("N21").Value = -81.002140249276
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 36
$ws.Range("K22").Value = -12.195121951219
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 71.428571428571
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 100
$ws.Range("L23").Value = -10.526315789473
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -21.052631578947
$ws.Range("F24").Value = 98
$ws.Range("H24").Value = -2.970297029702
$ws.Range("I24").Value = 1146
$ws.Range("J24").Value = 1361
$ws.Range("K24").Value = -15.797207935341
$ws.Range("L24").Value = -20.194986072423
$ws.Range("M24").Value = 21.52704135737
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -16.666666666666
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = -39.473684210526
$ws.Range("I25").Value = 469
$ws.Range("J25").Value = 505
$ws.Range("K25").Value = -7.128712871287
$ws.Range("L25").Value = 8.314087759815
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -22.222222222222
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = 12.903225806451
$ws.Range("I26").Value = 501
$ws.Range("J26").Value = 462
$ws.Range("K26").Value = 8.441558441558
$ws.Range("L26").Value = 15.437788018433
$ws.Range("M26").Value = 29.457364341085
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 8
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 43
$ws.Range("J28").Value = 48
$ws.Range("K28").Value = -10.416666666666
$ws.Range("L28").Value = -21.818181818181
$ws.Range("D31").Value = 2
$ws.Range("G31").Value = 4
$ws.Range("J31").Value = 28
$ws.Range("K31").Value = -75
$ws.Range("L33").Value = -28.571428571428
